{"js": "// Word JS API (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Applies the textual changes described by the commit \"Added many more\n// features\": retitles the review, rewrites the pros/cons bullet lists, and\n// rewrites the closing summary paragraph.\n\nconst replacements = [\n  [\n    \"Play Book of Ba Free - Exciting Egyptian-themed slot game\",\n    \"Play Book of Ba Free - Exciting Underworld Adventure Slot\",\n  ],\n  [\n    \"Exciting theme and graphics\",\n    \"Thrilling adventure to the underworld\",\n  ],\n  [\n    \"Special expanding symbols in the bonus round\",\n    \"Wilds, scatters, free spins, and special expanding symbols\",\n  ],\n  [\n    \"Shatter mode increases win multiplier up to 4x\",\n    \"Shatter mode creates new winning combinations\",\n  ],\n  [\n    \"Re-triggerable free spins feature\",\n    \"Immersive graphic design and adventurous background music\",\n  ],\n  [\n    \"Only 9 paylines may limit gameplay options\",\n    \"Limited number of paylines\",\n  ],\n  [\n    \"May not be appealing to those not interested in Egyptian mythology\",\n    \"No progressive jackpot\",\n  ],\n  [\n    \"Take a thrilling adventure to the underworld in this online slot game. Play free Book of Ba now and benefit from special expanding symbols and free spins.\",\n    \"Join the adventure and play Book of Ba for free. Try your luck with wilds, scatters, and free spins.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Applies the textual changes described by the commit \"Added many more\n# features\": retitles the review, rewrites the pros/cons bullet lists, and\n# rewrites the closing summary paragraph.\n#\n# $word / $d (ActiveDocument) are pre-bound by the host.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Book of Ba Free - Exciting Egyptian-themed slot game\", \"Play Book of Ba Free - Exciting Underworld Adventure Slot\"),\n    @(\"Exciting theme and graphics\", \"Thrilling adventure to the underworld\"),\n    @(\"Special expanding symbols in the bonus round\", \"Wilds, scatters, free spins, and special expanding symbols\"),\n    @(\"Shatter mode increases win multiplier up to 4x\", \"Shatter mode creates new winning combinations\"),\n    @(\"Re-triggerable free spins feature\", \"Immersive graphic design and adventurous background music\"),\n    @(\"Only 9 paylines may limit gameplay options\", \"Limited number of paylines\"),\n    @(\"May not be appealing to those not interested in Egyptian mythology\", \"No progressive jackpot\"),\n    @(\"Take a thrilling adventure to the underworld in this online slot game. Play free Book of Ba now and benefit from special expanding symbols and free spins.\", \"Join the adventure and play Book of Ba for free. Try your luck with wilds, scatters, and free spins.\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
